$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price (D) column cells that are being updated,
# so Excel does not auto-convert numeric-looking strings into real numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.345.43"
$ws.Range("E2").Value = "  -3.19%  "

$ws.Range("D3").Value = "1.937.17"
$ws.Range("E3").Value = "  -3.48%  "

$ws.Range("D4").Value = "1.000"

$ws.Range("D5").Value = "249.78"
$ws.Range("E5").Value = "  -3.88%  "

$ws.Range("D6").Value = "0.7241"
$ws.Range("E6").Value = "  -7.19%  "

$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "0.3338"
$ws.Range("E8").Value = "  -6.73%  "

$ws.Range("D9").Value = "28.55"
$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("D10").Value = "0.07460"
$ws.Range("E10").Value = "  +5.52%  "

$ws.Range("D11").Value = "0.8151"
$ws.Range("E11").Value = "  -5.14%  "

$ws.Range("D12").Value = "0.08128"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.512"
$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.932.25"
$ws.Range("E14").Value = "  -3.80%  "

$ws.Range("D15").Value = "95.05"
$ws.Range("E15").Value = "  -6.48%  "

$ws.Range("D16").Value = "14.85"
$ws.Range("E16").Value = "  -4.83%  "

$ws.Range("D17").Value = "30.345.35"
$ws.Range("E17").Value = "  -3.17%  "

$ws.Range("D18").Value = "0.000008380"
$ws.Range("E18").Value = "  +4.60%  "

$ws.Range("D19").Value = "254.74"
$ws.Range("E19").Value = "  -7.30%  "

$ws.Range("D20").Value = "5.872"
$ws.Range("E20").Value = "  -1.37%  "

$ws.Range("D21").Value = "2.189.20"
$ws.Range("E21").Value = "  -3.58%  "

$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "0.9996"

$ws.Range("D24").Value = "6.947"
$ws.Range("E24").Value = "  -3.25%  "

$ws.Range("D25").Value = "9.808"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "160.45"
$ws.Range("E26").Value = "  -3.37%  "

$ws.Range("D27").Value = "2.419"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").Value = "19.45"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").Value = "0.1334"
$ws.Range("E29").Value = "  -10.23%  "

$ws.Range("D30").Value = "1.561"
$ws.Range("E30").Value = "  -3.92%  "

$ws.Range("D31").Value = "1.343"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "4.444"
$ws.Range("E32").Value = "  -4.18%  "

$ws.Range("D33").Value = "4.238"
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("D34").Value = "0.05197"
$ws.Range("E34").Value = "  -0.77%  "

$ws.Range("D35").Value = "1.247"
$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("D36").Value = "0.7506"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("D37").Value = "2.734"
$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("D38").Value = "0.01998"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").Value = "2.838"
$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("D40").Value = "6.670"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").Value = "79.38"
$ws.Range("E41").Value = "  -0.96%  "

$ws.Range("D42").Value = "0.4536"
$ws.Range("E42").Value = "  -4.53%  "

$ws.Range("D43").Value = "2.026"
$ws.Range("E43").Value = "  -6.15%  "

$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").Value = "0.8399"
$ws.Range("E45").Value = "  -2.38%  "

$ws.Range("D46").Value = "102.59"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("D47").Value = "9.813"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Value = "7.387"
$ws.Range("E48").Value = "  -5.81%  "

$ws.Range("D49").Value = "36.82"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "1.506"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "0.4130"
$ws.Range("E51").Value = "  -5.45%  "

# Restore default (Normal) style on the Price column cells so no stray
# number-format style lingers on them, matching the original formatting.
foreach ($addr in $priceCells) { $ws.Range($addr).Style = "Normal" }
